# Auto commit at 2025-10-11 16:55:34.83
# Updates the "Metrics" sheet's monthly/yearly/total figures and fills in
# today's (daily) charge/income/order figures on the "today" sheet, then
# moves the active-tab/selection from "today" back to "Metrics".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metrics sheet: refresh the metric values in column B (rows 2-13)
# ---------------------------------------------------------------------
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 140784.92000000001
$metrics.Range("B3").Value  = 118549.43000000002
$metrics.Range("B4").Value  = 52007.28
$metrics.Range("B5").Value  = 5488
$metrics.Range("B6").Value  = 4507916.3900000006
$metrics.Range("B7").Value  = 3808368.0999999996
$metrics.Range("B8").Value  = 1322609.4200000002
$metrics.Range("B9").Value  = 174489
$metrics.Range("B10").Value = 32973240.190999825
$metrics.Range("B11").Value = 31083589.620000005
$metrics.Range("B12").Value = 11604318.309999999
$metrics.Range("B13").Value = 1272116

# ---------------------------------------------------------------------
# 2. today sheet: populate the day's figures (B3:B6), previously blank.
#    Dependent formulas (B11:B22, E11:E22, F11:F22, A1) recompute
#    automatically off of these + the Metrics values above.
# ---------------------------------------------------------------------
$today = $wb.Worksheets.Item("today")

$today.Range("B3").Value = 15411.33
$today.Range("B4").Value = 12836.54
$today.Range("B5").Value = 4934.28
$today.Range("B6").Value = 631

# Selection on "today" moves to F11:F22 and it is no longer the active tab.
$today.Range("F11:F22").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. Make "Metrics" the active sheet again, with D7 selected (was D8).
# ---------------------------------------------------------------------
$metrics.Activate() | Out-Null
$metrics.Range("D7").Select() | Out-Null
